# Add a new "correction" column (T) to the sheet, populate the data rows
# with the constant correction factor 0.782, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell T1, styled like the other header cells (centered).
$ws.Range("T1").Value = "correction"
$ws.Range("T1").HorizontalAlignment = -4108

# Fill T2:T9 with the correction value (plain numbers, no special style).
$ws.Range("T2:T9").Value = 0.782

# Update the active selection to match the final workbook state.
$ws.Range("R14").Select()
